$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that no longer belong in the fraud report ---
# (Spirometry, Throat culture, Upper arm X-ray, Intubation)
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("18").Delete()
$ws.Rows("16").Delete()
$ws.Rows("15").Delete()
$ws.Rows("13").Delete()

# After the deletes the remaining rows are:
#   9  Chlamydia antigen test
#   10 Allergy screening test
#   11 Intramuscular injection
#   12 Colonoscopy
#   13 Standard pregnancy test
#   14 Prostatectomy
#   15 Total Invoice Amount

# --- Insert a new line item before the total row ---
$ws.Rows("15").Insert()

# Match the "Legitimate" (green) formatting used elsewhere in the table
$ws.Range("A10:E10").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value = "seasonique 91 day pack"
$ws.Range("B15").Value = 500
$ws.Range("C15").Value = 431.4
$ws.Range("D15").Value = 68.59999999999999
$ws.Range("E15").Value = "Legitimate"

# --- Update the invoice total to reflect the new set of line items ---
$ws.Range("B16").Value = 47633.1
